$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits alone in the empty
#    paragraph right after the title (it will be re-created further
#    down, inside the "cobertura dos testes" paragraph).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Center-align the numeric values in the second column of the
#    summary table (32 / 26 / 4 / 22 / 6).
# ------------------------------------------------------------------
$t = $d.Tables(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cell = $t.Cell($i, 2)
    $cell.Range.Paragraphs(1).Alignment = 1
}

# ------------------------------------------------------------------
# 3. Rewrite the "cobertura dos testes" paragraph with the updated
#    figures / wording. This is done in two pieces so the existing
#    <w:proofErr/> marks (around the double space that follows
#    "trinta"/"seis") are left untouched, just like in the target.
# ------------------------------------------------------------------
$oldA = "sendo que trinta"
$newA = "sendo que foram executados Vinte e seis"
$d.Content.Find.Execute($oldA, $true, $false, $false, $false, $false, $true, 1, $false, $newA, 2)

$oldB = "( 30  ) casos de testes, sendo oito(08) testes com sucesso onze(10) com falha e quatro (2) não foram executados, "
$newB = "(26 ) casos de testes, sendo quatro (04) testes com sucesso, Vinte e dois (22) com falha e seis (06) não foram executados, "
$d.Content.Find.Execute($oldB, $true, $false, $false, $false, $false, $true, 1, $false, $newB, 2)

# ------------------------------------------------------------------
# 4. Re-insert the "_GoBack" bookmark so it starts right after the
#    paragraph's leading tab and ends at the end of the paragraph
#    (this is where Word had left the editing cursor).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Os testes foram satisfatórios, sendo que foram executados") | Out-Null
$bmStart = $rng.Start

$rng2 = $d.Content
$rng2.Find.Execute("pois o programa só estava com a parte gerencial.") | Out-Null
$bmEnd = $rng2.End

$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
